# repull data, push all data, mean calculation
# Update the dSF (column F) values on Sheet1 to reflect the repulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "F3"  = 1
    "F5"  = 1
    "F6"  = -1
    "F8"  = 5
    "F10" = -7
    "F11" = -2
    "F12" = 0
    "F15" = 0
    "F18" = 2
    "F19" = -4
    "F20" = -5
    "F21" = -2
    "F25" = 0
    "F27" = -3
    "F28" = 8
    "F30" = 7
    "F32" = -6
    "F34" = -3
    "F37" = -2
    "F38" = 9
    "F39" = 0
    "F40" = 2
    "F42" = -7
    "F43" = 2
    "F44" = 6
    "F47" = 2
    "F48" = 4
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
